$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translatable_Site_labels")

# Insert 4 new rows right after the current "undoEditBtn" row (rows 57..60),
# pushing the trailing "}" / "export default labels;" rows down to 61 / 63
# (keeping the existing blank-row gap before the last line).
$ws.Rows("57:60").Insert()

# New "Site Header" group of translatable labels (Label id / Description / group).
# Column D (group/description) first, then B (ids), then C (display text) -- this
# mirrors the order the original author typed them in and keeps the resulting
# shared-string table ordering aligned with the source workbook.
$ws.Range("D57").Value2 = "Site Header"
$ws.Range("D58").Value2 = "Site Header"
$ws.Range("D59").Value2 = "Site Header"
$ws.Range("D60").Value2 = "Site Header"

$ws.Range("B57").Value2 = "homeBtn"
$ws.Range("B58").Value2 = "detailedSearchBtn"
$ws.Range("B59").Value2 = "importDataBtn"
$ws.Range("B60").Value2 = "adminBtn"

$ws.Range("C57").Value2 = "Home"
$ws.Range("C58").Value2 = "Detailed Search"
$ws.Range("C59").Value2 = "Import"
$ws.Range("C60").Value2 = "Admin"

# Re-create the generator formula (CONCAT of id/label pair) for the new rows --
# same pattern used throughout column E.
$ws.Range("E57:E60").Formula = "=_xlfn.CONCAT("""",B57,"" : '"",C57,""',"")"

# Move the selection/view to the top of the sheet.
$ws.Range("E10").Select() | Out-Null
